$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.996.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3
$ws.Range("D3").Value = "'1.851.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "

# Row 4
$ws.Range("D4").Value = "'1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "

# Row 5
$ws.Range("D5").Value = "'1.012"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "

# Row 6
$ws.Range("D6").Value = "'309.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

# Row 7
$ws.Range("D7").Value = "'0.4770"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.88%  "

# Row 8
$ws.Range("D8").Value = "'0.3676"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.01%  "

# Row 9
$ws.Range("D9").Value = "'0.07228"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.27%  "

# Row 10
$ws.Range("D10").Value = "'0.9302"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.97%  "

# Row 11
$ws.Range("D11").Value = "'19.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.54%  "

# Row 12
$ws.Range("D12").Value = "'0.07729"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "

# Row 13
$ws.Range("D13").Value = "'1.898.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.51%  "

# Row 14
$ws.Range("D14").Value = "'5.327"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15
$ws.Range("D15").Value = "'6.422"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.18%  "

# Row 16
$ws.Range("D16").Value = "'88.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "

# Row 17
$ws.Range("D17").Value = "'1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "

# Row 18
$ws.Range("D18").Value = "'0.000008644"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "

# Row 19
$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20
$ws.Range("D20").Value = "'27.023.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "

# Row 21
$ws.Range("D21").Value = "'14.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "

# Row 22
$ws.Range("D22").Value = "'5.062"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "

# Row 23
$ws.Range("D23").Value = "'10.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "

# Row 24
$ws.Range("D24").Value = "'1.937"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "

# Row 25
$ws.Range("D25").Value = "'152.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "

# Row 26
$ws.Range("D26").Value = "'18.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("D27").Value = "'2.003"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.58%  "

# Row 28
$ws.Range("D28").Value = "'114.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.31%  "

# Row 29
$ws.Range("D29").Value = "'4.983"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.89%  "

# Row 30
$ws.Range("D30").Value = "'0.08888"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.84%  "

# Row 31
$ws.Range("D31").Value = "'3.328"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.70%  "

# Row 32
$ws.Range("D32").Value = "'1.173"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "

# Row 33
$ws.Range("D33").Value = "'0.7453"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.16%  "

# Row 34
$ws.Range("D34").Value = "'4.504"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.20%  "

# Row 35
$ws.Range("D35").Value = "'2.732"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.26%  "

# Row 36
$ws.Range("D36").Value = "'1.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.89%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05280"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.50%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01956"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.36%  "

# Row 39
$ws.Range("D39").Value = "'2.972"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.54%  "

# Row 40
$ws.Range("D40").Value = "'0.5194"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "

# Row 41
$ws.Range("D41").Value = "'7.013"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.94%  "

# Row 42
$ws.Range("D42").Value = "'0.1511"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "

# Row 43
$ws.Range("D43").Value = "'8.208"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "

# Row 44
$ws.Range("D44").Value = "'10.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.83%  "

# Row 45
$ws.Range("D45").Value = "'0.4743"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "

# Row 46
$ws.Range("D46").Value = "'1.014"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.53%  "

# Row 47
$ws.Range("D47").Value = "'101.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.61%  "

# Row 48
$ws.Range("D48").Value = "'1.607"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.42%  "

# Row 49
$ws.Range("D49").Value = "'65.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.50%  "

# Row 50
$ws.Range("D50").Value = "'0.06029"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "

# Row 51
$ws.Range("D51").Value = "'0.8886"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.30%  "
